$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet/tab to reflect the new "through" date
$ws.Name = "Through 2021-09-29"

# Row 9 (July) - 2021 arrest-made/no-arrest-made/rate columns updated
$ws.Range("T9").Value = 11
$ws.Range("U9").Value = 140
$ws.Range("V9").Value = 0.0728

# Row 11 (September) label updated to reflect new "through" date
$ws.Range("A11").Value = "September (through 09-29)"

# Row 11 (September) data updates across all years
$ws.Range("C11").Value = 29
$ws.Range("D11").Value = 0.0333
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 39
$ws.Range("G11").Value = 0.1136
$ws.Range("I11").Value = 69
$ws.Range("J11").Value = 0.0676
$ws.Range("L11").Value = 50
$ws.Range("M11").Value = 0.0741
$ws.Range("O11").Value = 64
$ws.Range("P11").Value = 0.0986
$ws.Range("Q11").Value = 5
$ws.Range("R11").Value = 107
$ws.Range("S11").Value = 0.0446
$ws.Range("U11").Value = 170
$ws.Range("V11").Value = 0.0116

# Row 12 (Total) data updates across all years
$ws.Range("C12").Value = 194
$ws.Range("D12").Value = 0.1339
$ws.Range("E12").Value = 46
$ws.Range("F12").Value = 379
$ws.Range("G12").Value = 0.1082
$ws.Range("I12").Value = 575
$ws.Range("J12").Value = 0.08
$ws.Range("L12").Value = 483
$ws.Range("M12").Value = 0.1121
$ws.Range("O12").Value = 377
$ws.Range("P12").Value = 0.1024
$ws.Range("Q12").Value = 53
$ws.Range("R12").Value = 843
$ws.Range("S12").Value = 0.0592
$ws.Range("T12").Value = 78
$ws.Range("U12").Value = 1164
$ws.Range("V12").Value = 0.0628
